$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above row 622, shifting the existing rows
# (622-705) down to (625-708).
$ws.Range("A622:A624").EntireRow.Insert()

# Populate the 3 newly inserted rows with the new data records.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44776, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Primera", 108, 20000, 20000, 20000, "$/caja 12 unidades", "Ecuador", 1667, 12),
    @(3, "Femacal de La Calera", "Coquimbo", 44776, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Segunda", 108, 20000, 20000, 20000, "$/caja 14 unidades", "Ecuador", 1429, 14),
    @(3, "Femacal de La Calera", "Coquimbo", 44776, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108005, "Piña", "Caramelo", "Tercera", 54,  20000, 20000, 20000, "$/caja 16 unidades", "Ecuador", 1250, 16)
)

for ($i = 0; $i -lt 3; $i++) {
    $r = 622 + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
